$d = $word.ActiveDocument

# --- 1. Insert a new paragraph right after the final ("(SWE 1)") paragraph ---
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()

# The new paragraph inherits the same pPr/rPr (style "Paragraphedeliste",
# ind left=810, italic) as the paragraph it was split from.
$newPara = $d.Paragraphs.Last

# A sentinel character is appended after the real text; this keeps the
# "after text" insertion point for the bookmark away from the very end of
# the document's content while we create it (doing it exactly at end-of-
# story  can mis-place the bookmark at the very start of the document).
# The sentinel is stripped again right afterwards.
$newPara.Range.Text = "MATRICLE:ET20210153X"

# --- 2. Move the _GoBack bookmark from the old last paragraph onto the new one ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$newPara = $d.Paragraphs.Last
$bmRange = $newPara.Range.Duplicate
[void]$bmRange.MoveEnd(1, -1)   # exclude the paragraph mark
[void]$bmRange.MoveEnd(1, -1)   # exclude the trailing sentinel character
$bmRange.Collapse(0)            # wdCollapseEnd: collapse to a point right after the real text
$d.Bookmarks.Add("_GoBack", $bmRange)

# --- 3. Strip the sentinel character back out ---
$newPara = $d.Paragraphs.Last
$full = $newPara.Range.Duplicate
[void]$full.MoveEnd(1, -1)      # exclude the paragraph mark
$sentinel = $d.Range($full.End - 1, $full.End)
$sentinel.Text = ""
